$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 939.0952
$ws.Range("J17").Value = 868
$ws.Range("L17").Value = 2604
$ws.Range("N17").Value = -2940
$ws.Range("H58").Value = 538.1111
$ws.Range("I58").Value = 57.333332
$ws.Range("J58").Value = 1499.6666
$ws.Range("K58").Value = 171.999996
$ws.Range("L58").Value = 4498.9998
$ws.Range("M58").Value = -21.99999600000001
$ws.Range("N58").Value = -4798.9998
$ws.Range("H86").Value = 2500
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 2500
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 2500
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -4746
$ws.Range("H89").Value = 2500
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 2500
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 12500
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -23732
$ws.Range("H136").Value = 80000
$ws.Range("J136").Value = 80000
$ws.Range("L136").Value = 80000
$ws.Range("N136").Value = -90200
$ws.Range("H137").Value = 2655.024
$ws.Range("I137").Value = 2197.7368
$ws.Range("K137").Value = 6593.2104
$ws.Range("M137").Value = -4043.2104
$ws.Range("H138").Value = 2847.3635
$ws.Range("I138").Value = 1477.9166
$ws.Range("J138").Value = 3629.9048
$ws.Range("K138").Value = 4433.7498
$ws.Range("L138").Value = 10889.7144
$ws.Range("M138").Value = 706.2502000000004
$ws.Range("N138").Value = -21169.7144

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 24394646
$ws.Range("I32").Value = 26320262
$ws.Range("K32").Value = 26320262
$ws.Range("M32").Value = -26319975
$ws.Range("H61").Value = 3274.375
$ws.Range("I61").Value = 3274.375
$ws.Range("K61").Value = 3274.375
$ws.Range("M61").Value = -3062.375
$ws.Range("H97").Value = 895.3333
$ws.Range("I97").Value = 940.86664
$ws.Range("K97").Value = 940.86664
$ws.Range("M97").Value = -444.86664
$ws.Range("H132").Value = 1899.6875
$ws.Range("I132").Value = 1899.6875
$ws.Range("K132").Value = 5699.0625
$ws.Range("M132").Value = -3169.0625
$ws.Range("H136").Value = 3274.375
$ws.Range("I136").Value = 3274.375
$ws.Range("K136").Value = 9823.125
$ws.Range("M136").Value = -7273.125

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 3268.6667
$ws.Range("I94").Value = 2930.8572
$ws.Range("J94").Value = 3944.2856
$ws.Range("K94").Value = 2930.8572
$ws.Range("L94").Value = 3944.2856
$ws.Range("M94").Value = -2479.8572
$ws.Range("N94").Value = -4846.2856
$ws.Range("H105").Value = 1741.1923
$ws.Range("I105").Value = 1765.9412
$ws.Range("J105").Value = 1694.4445
$ws.Range("K105").Value = 1765.9412
$ws.Range("L105").Value = 1694.4445
$ws.Range("M105").Value = -18.94119999999998
$ws.Range("N105").Value = -5188.4445
$ws.Range("H134").Value = 2780.3635
$ws.Range("I134").Value = 2119.3684
$ws.Range("J134").Value = 6966.6665
$ws.Range("K134").Value = 6358.1052
$ws.Range("L134").Value = 20899.9995
$ws.Range("M134").Value = -3823.1052
$ws.Range("N134").Value = -25969.9995

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 7562.8667
$ws.Range("I16").Value = 8453.75
$ws.Range("K16").Value = 8453.75
$ws.Range("M16").Value = -8166.75
$ws.Range("H31").Value = 1382.619
$ws.Range("I31").Value = 1251.75
$ws.Range("J31").Value = 4000
$ws.Range("K31").Value = 1251.75
$ws.Range("L31").Value = 4000
$ws.Range("M31").Value = -956.75
$ws.Range("N31").Value = -4590
$ws.Range("H34").Value = 1382.619
$ws.Range("I34").Value = 1251.75
$ws.Range("J34").Value = 4000
$ws.Range("K34").Value = 1251.75
$ws.Range("L34").Value = 4000
$ws.Range("M34").Value = -1049.75
$ws.Range("N34").Value = -4404
$ws.Range("H58").Value = 2751.5386
$ws.Range("I58").Value = 1918.8889
$ws.Range("J58").Value = 4625
$ws.Range("K58").Value = 1918.8889
$ws.Range("L58").Value = 4625
$ws.Range("M58").Value = -1715.8889
$ws.Range("N58").Value = -5031
$ws.Range("H95").Value = 22500
$ws.Range("J95").Value = 22500
$ws.Range("L95").Value = 22500
$ws.Range("N95").Value = -27992
$ws.Range("H99").Value = 14963735
$ws.Range("I99").Value = 2714558
$ws.Range("K99").Value = 2714558
$ws.Range("M99").Value = -2713060
$ws.Range("H113").Value = 7562.8667
$ws.Range("I113").Value = 8453.75
$ws.Range("K113").Value = 8453.75
$ws.Range("M113").Value = -6283.75
$ws.Range("H126").Value = 14963735
$ws.Range("I126").Value = 2714558
$ws.Range("K126").Value = 8143674
$ws.Range("M126").Value = -8141204
$ws.Range("H132").Value = 2862.138
$ws.Range("I132").Value = 2265.5652
$ws.Range("K132").Value = 6796.6956
$ws.Range("M132").Value = -4266.6956
$ws.Range("H134").Value = 7907.615
$ws.Range("I134").Value = 7907.615
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 23722.845
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -21187.845
$ws.Range("N134").ClearContents()
$ws.Range("H136").Value = 2751.5386
$ws.Range("I136").Value = 1918.8889
$ws.Range("J136").Value = 4625
$ws.Range("K136").Value = 5756.6667
$ws.Range("L136").Value = 13875
$ws.Range("M136").Value = -3206.6667
$ws.Range("N136").Value = -18975

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 99
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()
$ws.Range("H7").Value = 7743.4287
$ws.Range("J7").Value = 942.8570999999999
$ws.Range("L7").Value = 2828.5713
$ws.Range("N7").Value = -3052.5713
$ws.Range("H8").Value = 1675.2222
$ws.Range("I8").Value = 1675.2222
$ws.Range("K8").Value = 5025.6666
$ws.Range("M8").Value = -4886.6666
$ws.Range("H117").Value = 1521.091
$ws.Range("I117").Value = 400
$ws.Range("J117").Value = 2455.3333
$ws.Range("K117").Value = 1200
$ws.Range("L117").Value = 7365.999899999999
$ws.Range("M117").Value = 2242
$ws.Range("N117").Value = -14249.9999
$ws.Range("H140").Value = 5545.1816
$ws.Range("I140").Value = 2599.4
$ws.Range("K140").Value = 7798.200000000001
$ws.Range("M140").Value = -2618.200000000001

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 1263.375
$ws.Range("J10").Value = 767.8333
$ws.Range("L10").Value = 767.8333
$ws.Range("N10").Value = -1105.8333
$ws.Range("H19").Value = 2500376
$ws.Range("I19").Value = 3333751.2
$ws.Range("J19").Value = 250
$ws.Range("K19").Value = 3333751.2
$ws.Range("L19").Value = 250
$ws.Range("M19").Value = -3333463.2
$ws.Range("N19").Value = -826
$ws.Range("H44").Value = 25499
$ws.Range("J44").Value = 25499
$ws.Range("L44").Value = 25499
$ws.Range("N44").Value = -26691
$ws.Range("H102").Value = 3505.0344
$ws.Range("I102").Value = 2529.5454
$ws.Range("K102").Value = 2529.5454
$ws.Range("M102").Value = -907.5454
$ws.Range("H122").Value = 4821.7144
$ws.Range("I122").Value = 4117.846
$ws.Range("K122").Value = 12353.538
$ws.Range("M122").Value = -9903.537999999999
$ws.Range("H126").Value = 5343
$ws.Range("I126").Value = 4931
$ws.Range("K126").Value = 14793
$ws.Range("M126").Value = -12323
$ws.Range("H132").Value = 4131.5835
$ws.Range("I132").Value = 4511.05
$ws.Range("J132").Value = 2234.25
$ws.Range("K132").Value = 13533.15
$ws.Range("L132").Value = 6702.75
$ws.Range("M132").Value = -11003.15
$ws.Range("N132").Value = -11762.75

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H47").Value = 38165
$ws.Range("J47").Value = 38165
$ws.Range("L47").Value = 38165
$ws.Range("N47").Value = -39145
$ws.Range("H52").Value = 38165
$ws.Range("J52").Value = 38165
$ws.Range("L52").Value = 38165
$ws.Range("N52").Value = -38631
$ws.Range("H93").Value = 9007.875
$ws.Range("I93").Value = 8149.25
$ws.Range("J93").Value = 9866.5
$ws.Range("K93").Value = 8149.25
$ws.Range("L93").Value = 9866.5
$ws.Range("M93").Value = -6901.25
$ws.Range("N93").Value = -12362.5
$ws.Range("H132").Value = 2213.0425
$ws.Range("I132").Value = 1813.1666
$ws.Range("K132").Value = 5439.4998
$ws.Range("M132").Value = -2909.4998

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 889
$ws.Range("I107").Value = 847.7778
$ws.Range("J107").Value = 942
$ws.Range("K107").Value = 2543.3334
$ws.Range("L107").Value = 2826
$ws.Range("M107").Value = -623.3334
$ws.Range("N107").Value = -6666
$ws.Range("H126").Value = 2509.818
$ws.Range("I126").Value = 2622.5557
$ws.Range("J126").Value = 2002.5
$ws.Range("K126").Value = 7867.6671
$ws.Range("L126").Value = 6007.5
$ws.Range("M126").Value = -5397.6671
$ws.Range("N126").Value = -10947.5
$ws.Range("H132").Value = 914.54285
$ws.Range("I132").Value = 914.54285
$ws.Range("K132").Value = 2743.62855
$ws.Range("M132").Value = -213.6285500000004
$ws.Range("H135").Value = 61116.824
$ws.Range("J135").Value = 61116.824
$ws.Range("L135").Value = 61116.824
$ws.Range("N135").Value = -71256.82399999999
